$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 30.75612566666667
    "H2" = 92.268377
    "I2" = 0.9777985798685588
    "J2" = 0.9777985798685588
    "M2" = 2.231113333333334
    "N2" = 6.69334
    "O2" = 0.01598125358798882
    "P2" = 0.01598125358798882
    "Q2" = 68.62040205657556
    "R2" = 617.58361850918
    "S2" = 0.01562644706285478
    "T2" = 0.01562644706285478
    "G3" = 30.75612566666667
    "H3" = 92.268377
    "I3" = 0.9777985798685588
    "J3" = 0.9777985798685588
    "O3" = 0.1634493267640196
    "P3" = 0.1634493267640195
    "Q3" = 701.8196949739477
    "R3" = 6316.377254765529
    "S3" = 0.1598205195903304
    "T3" = 0.1598205195903303
    "G4" = 30.75612566666667
    "H4" = 92.268377
    "I4" = 0.9777985798685588
    "J4" = 0.9777985798685588
    "M4" = 58.02175166666666
    "N4" = 174.065255
    "O4" = 0.4156043142904646
    "P4" = 0.4156043142904646
    "Q4" = 1784.524285660126
    "R4" = 16060.71857094113
    "S4" = 0.4063773083004624
    "T4" = 0.4063773083004624
    "G5" = 30.75612566666667
    "H5" = 92.268377
    "I5" = 0.9777985798685588
    "J5" = 0.9777985798685588
    "M5" = 15.16934033333333
    "N5" = 45.508021
    "O5" = 0.1086565487318021
    "P5" = 0.1086565487318021
    "Q5" = 466.5501375724352
    "R5" = 4198.951238151917
    "S5" = 0.106244219043375
    "T5" = 0.106244219043375
    "G6" = 30.75612566666667
    "H6" = 92.268377
    "I6" = 0.9777985798685588
    "J6" = 0.9777985798685588
    "M6" = 41.36709099999999
    "N6" = 124.101273
    "O6" = 0.2963085566257249
    "P6" = 0.2963085566257249
    "Q6" = 1272.291449260435
    "R6" = 11450.62304334392
    "S6" = 0.2897300858715363
    "T6" = 0.2897300858715363
    "I7" = 0.004830327290741966
    "J7" = 0.004830327290741966
    "M7" = 2.231113333333334
    "N7" = 6.69334
    "O7" = 0.01598125358798882
    "P7" = 0.01598125358798882
    "Q7" = 0.3389849480044445
    "R7" = 3.05086453204
    "S7" = 0.00007719468534633035
    "T7" = 0.00007719468534633035
    "I8" = 0.004830327290741966
    "J8" = 0.004830327290741966
    "O8" = 0.1634493267640196
    "P8" = 0.1634493267640195
    "S8" = 0.000789513743721645
    "T8" = 0.0007895137437216449
    "I9" = 0.004830327290741966
    "J9" = 0.004830327290741966
    "M9" = 58.02175166666666
    "N9" = 174.065255
    "O9" = 0.4156043142904646
    "P9" = 0.4156043142904646
    "Q9" = 8.815554180058887
    "R9" = 79.33998762053
    "S9" = 0.002007504861467332
    "T9" = 0.002007504861467332
    "I10" = 0.004830327290741966
    "J10" = 0.004830327290741966
    "M10" = 15.16934033333333
    "N10" = 45.508021
    "O10" = 0.1086565487318021
    "P10" = 0.1086565487318021
    "Q10" = 2.304758779991778
    "R10" = 20.742829019926
    "S10" = 0.0005248466926570582
    "T10" = 0.0005248466926570582
    "I11" = 0.004830327290741966
    "J11" = 0.004830327290741966
    "M11" = 41.36709099999999
    "N11" = 124.101273
    "O11" = 0.2963085566257249
    "P11" = 0.2963085566257249
    "Q11" = 6.285122760115333
    "R11" = 56.566104841038
    "S11" = 0.0014312673075496
    "T11" = 0.0014312673075496
    "E12" = 3
    "F12" = 1
    "G12" = 0.4673666666666667
    "H12" = 1.4021
    "I12" = 0.01485851852399773
    "J12" = 0.01485851852399773
    "M12" = 2.231113333333334
    "N12" = 6.69334
    "O12" = 0.01598125358798882
    "P12" = 0.01598125358798882
    "Q12" = 1.042748001555556
    "R12" = 9.384732014000001
    "S12" = 0.0002374577524738371
    "T12" = 0.0002374577524738371
    "E13" = 3
    "F13" = 1
    "G13" = 0.4673666666666667
    "H13" = 1.4021
    "I13" = 0.01485851852399773
    "J13" = 0.01485851852399773
    "O13" = 0.1634493267640196
    "P13" = 0.1634493267640195
    "Q13" = 10.66477406796667
    "R13" = 95.98296661170001
    "S13" = 0.002428614849458144
    "T13" = 0.002428614849458143
    "E14" = 3
    "F14" = 1
    "G14" = 0.4673666666666667
    "H14" = 1.4021
    "I14" = 0.01485851852399773
    "J14" = 0.01485851852399773
    "M14" = 58.02175166666666
    "N14" = 174.065255
    "O14" = 0.4156043142904646
    "P14" = 0.4156043142904646
    "Q14" = 27.11743267061111
    "R14" = 244.0568940355
    "S14" = 0.006175264402538244
    "T14" = 0.006175264402538244
    "E15" = 3
    "F15" = 1
    "G15" = 0.4673666666666667
    "H15" = 1.4021
    "I15" = 0.01485851852399773
    "J15" = 0.01485851852399773
    "M15" = 15.16934033333333
    "N15" = 45.508021
    "O15" = 0.1086565487318021
    "P15" = 0.1086565487318021
    "Q15" = 7.089644027122223
    "R15" = 63.80679624410001
    "S15" = 0.001614475342085144
    "T15" = 0.001614475342085144
    "E16" = 3
    "F16" = 1
    "G16" = 0.4673666666666667
    "H16" = 1.4021
    "I16" = 0.01485851852399773
    "J16" = 0.01485851852399773
    "M16" = 41.36709099999999
    "N16" = 124.101273
    "O16" = 0.2963085566257249
    "P16" = 0.2963085566257249
    "Q16" = 19.33359943036666
    "R16" = 174.0023948733
    "S16" = 0.004402706177442365
    "T16" = 0.004402706177442365
    "E17" = 2
    "F17" = 0.6666666666666666
    "G17" = 0.07903166666666667
    "H17" = 0.237095
    "I17" = 0.002512574316701549
    "J17" = 0.002512574316701549
    "M17" = 2.231113333333334
    "N17" = 6.69334
    "O17" = 0.01598125358798882
    "P17" = 0.01598125358798882
    "Q17" = 0.1763286052555556
    "R17" = 1.5869574473
    "S17" = 0.00004015408731387519
    "T17" = 0.00004015408731387519
    "E18" = 2
    "F18" = 0.6666666666666666
    "G18" = 0.07903166666666667
    "H18" = 0.237095
    "I18" = 0.002512574316701549
    "J18" = 0.002512574316701549
    "O18" = 0.1634493267640196
    "P18" = 0.1634493267640195
    "Q18" = 1.803412458201667
    "R18" = 16.230712123815
    "S18" = 0.0004106785805094347
    "T18" = 0.0004106785805094347
    "E19" = 2
    "F19" = 0.6666666666666666
    "G19" = 0.07903166666666667
    "H19" = 0.237095
    "I19" = 0.002512574316701549
    "J19" = 0.002512574316701549
    "M19" = 58.02175166666666
    "N19" = 174.065255
    "O19" = 0.4156043142904646
    "P19" = 0.4156043142904646
    "Q19" = 4.585555737136111
    "R19" = 41.270001634225
    "S19" = 0.00104423672599658
    "T19" = 0.00104423672599658
    "E20" = 2
    "F20" = 0.6666666666666666
    "G20" = 0.07903166666666667
    "H20" = 0.237095
    "I20" = 0.002512574316701549
    "J20" = 0.002512574316701549
    "M20" = 15.16934033333333
    "N20" = 45.508021
    "O20" = 0.1086565487318021
    "P20" = 0.1086565487318021
    "Q20" = 1.198858248777222
    "R20" = 10.789724238995
    "S20" = 0.0002730076536849563
    "T20" = 0.0002730076536849563
    "E21" = 2
    "F21" = 0.6666666666666666
    "G21" = 0.07903166666666667
    "H21" = 0.237095
    "I21" = 0.002512574316701549
    "J21" = 0.002512574316701549
    "M21" = 41.36709099999999
    "N21" = 124.101273
    "O21" = 0.2963085566257249
    "P21" = 0.2963085566257249
    "Q21" = 3.269310146881666
    "R21" = 29.423791321935
    "S21" = 0.0007444972691967031
    "T21" = 0.0007444972691967031
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
